# "Generate Report for Handback" -- mark the zh-cn / de-de localization rows
# as handed back (in sync with en-US), stamp the handback file name + the
# handback datetime, and widen the columns that now hold the longer text.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

$urlA = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/dbe4923a2fda080e891f7656ff1b1a97dcbafd53/e2e/a.md"
$urlB = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/dbe4923a2fda080e891f7656ff1b1a97dcbafd53/e2e/b.md"

# Hyperlink font color used by the workbook's "HyperLink" cell style
# (RGB 6495ED stored as a BGR COM color value).
$hyperlinkColor = 15570276

# ---------------------------------------------------------------------
# Overview sheet: just the Status columns (E = zh-cn, F = de-de) change
# text (via the shared "Status" string), which widens those columns.
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $statusText
$overview.Range("F2").Value = $statusText
$overview.Range("E3").Value = $statusText
$overview.Range("F3").Value = $statusText

$overview.Columns.Item(5).ColumnWidth = 29.144371396019366
$overview.Columns.Item(6).ColumnWidth = 29.144371396019366

# ---------------------------------------------------------------------
# zh-cn sheet: status -> handed back, target/handback file + datetime
# filled in for both data rows.
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $statusText
$zhcn.Range("C3").Value = $statusText

$zhcn.Range("I2").Value = "a.md"
$zhcn.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-23 16:40:05"

$zhcn.Range("I3").Value = "a.md"
$zhcn.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-08-23 16:40:05"

$zhcn.Hyperlinks.Delete()
$zhcn.Hyperlinks.Add($zhcn.Range("A2"), $urlA, [System.Type]::Missing, [System.Type]::Missing, "a.md")
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $urlA, [System.Type]::Missing, [System.Type]::Missing, "a.md")
$zhcn.Hyperlinks.Add($zhcn.Range("A3"), $urlB, [System.Type]::Missing, [System.Type]::Missing, "b.md")
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), $urlA, [System.Type]::Missing, [System.Type]::Missing, "a.md")
foreach ($addr in @("A2", "I2", "A3", "I3")) {
    $zhcn.Range($addr).Font.Underline = $true
    $zhcn.Range($addr).Font.Color = $hyperlinkColor
}

$zhcn.Columns.Item(3).ColumnWidth = 29.144371396019366
$zhcn.Columns.Item(10).ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------
# de-de sheet: same shape of update, using the de-de handback file name
# and its own handback datetime.
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $statusText
$dede.Range("C3").Value = $statusText

$dede.Range("I2").Value = "a.md"
$dede.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$dede.Range("K2").Value = "2016-08-23 16:40:30"

$dede.Range("I3").Value = "a.md"
$dede.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$dede.Range("K3").Value = "2016-08-23 16:40:30"

$dede.Hyperlinks.Delete()
$dede.Hyperlinks.Add($dede.Range("A2"), $urlA, [System.Type]::Missing, [System.Type]::Missing, "a.md")
$dede.Hyperlinks.Add($dede.Range("I2"), $urlA, [System.Type]::Missing, [System.Type]::Missing, "a.md")
$dede.Hyperlinks.Add($dede.Range("A3"), $urlB, [System.Type]::Missing, [System.Type]::Missing, "b.md")
$dede.Hyperlinks.Add($dede.Range("I3"), $urlA, [System.Type]::Missing, [System.Type]::Missing, "a.md")
foreach ($addr in @("A2", "I2", "A3", "I3")) {
    $dede.Range($addr).Font.Underline = $true
    $dede.Range($addr).Font.Color = $hyperlinkColor
}

$dede.Columns.Item(3).ColumnWidth = 29.144371396019366
$dede.Columns.Item(10).ColumnWidth = 39.166666666666664
